# Fruta / hortaliza, semanal
# Insert two new price records (rows 129-130) into the daily logic subset
# sheet for "Feria Lagunitas de Puerto Montt - Palta". Existing rows 129..205
# shift down to 131..207.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 129 - everything from the old
# row 129 onward shifts down by two rows (205 -> 207).
$ws.Rows("129:130").Insert()

# --- New row 129 ---
$ws.Cells.Item(129, 1).Value = 4
$ws.Cells.Item(129, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(129, 3).Value = "Los Lagos"
$ws.Cells.Item(129, 4).Value = 44460
$ws.Cells.Item(129, 5).Value = 10
$ws.Cells.Item(129, 6).Value = "Fruta"
$ws.Cells.Item(129, 7).Value = 100106
$ws.Cells.Item(129, 8).Value = "Oleaginosos"
$ws.Cells.Item(129, 9).Value = 100106002
$ws.Cells.Item(129, 10).Value = "Palta"
$ws.Cells.Item(129, 11).Value = "Hass"
$ws.Cells.Item(129, 12).Value = "1a nueva(o)"
$ws.Cells.Item(129, 13).Value = 80
$ws.Cells.Item(129, 14).Value = 4000
$ws.Cells.Item(129, 15).Value = 4000
$ws.Cells.Item(129, 16).Value = 4000
$ws.Cells.Item(129, 17).Value = "$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(129, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(129, 19).Value = 4000
$ws.Cells.Item(129, 20).Value = 1

# --- New row 130 ---
$ws.Cells.Item(130, 1).Value = 4
$ws.Cells.Item(130, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(130, 3).Value = "Los Lagos"
$ws.Cells.Item(130, 4).Value = 44460
$ws.Cells.Item(130, 5).Value = 10
$ws.Cells.Item(130, 6).Value = "Fruta"
$ws.Cells.Item(130, 7).Value = 100106
$ws.Cells.Item(130, 8).Value = "Oleaginosos"
$ws.Cells.Item(130, 9).Value = 100106002
$ws.Cells.Item(130, 10).Value = "Palta"
$ws.Cells.Item(130, 11).Value = "Hass"
$ws.Cells.Item(130, 12).Value = "2a nueva(o)"
$ws.Cells.Item(130, 13).Value = 70
$ws.Cells.Item(130, 14).Value = 3600
$ws.Cells.Item(130, 15).Value = 3600
$ws.Cells.Item(130, 16).Value = 3600
$ws.Cells.Item(130, 17).Value = "$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(130, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(130, 19).Value = 3600
$ws.Cells.Item(130, 20).Value = 1
